# Rolling-forward weekly forecast data one week, and refreshing the
# Summary sheet statistics to match ("Penalty Reward System" update).

$wb  = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Sheet "Forecast Comparison": update Week_Start_Date (B) and MyForecast (D)
# for rows 2-17. Dates are stored as plain text, so a leading apostrophe is
# used to stop Excel from auto-converting the date-like strings into date
# serial numbers.

$forecastUpdates = @{
    2  = @{ B = "2025-01-12"; D = 55 }
    3  = @{ B = "2025-01-19"; D = 43 }
    4  = @{ B = "2025-01-26"; D = 40 }
    5  = @{ B = "2025-02-02"; D = 41 }
    6  = @{ B = "2025-02-09"; D = 41 }
    7  = @{ B = "2025-02-16"; D = 42 }
    8  = @{ B = "2025-02-23"; D = 45 }
    9  = @{ B = "2025-03-02"; D = 34 }
    10 = @{ B = "2025-03-09"; D = 32 }
    11 = @{ B = "2025-03-16"; D = 33 }
    12 = @{ B = "2025-03-23"; D = 47 }
    13 = @{ B = "2025-03-30"; D = 48 }
    14 = @{ B = "2025-04-06"; D = 33 }
    15 = @{ B = "2025-04-13"; D = 32 }
    16 = @{ B = "2025-04-20"; D = 32 }
    17 = @{ B = "2025-04-27"; D = 31 }
}

foreach ($row in $forecastUpdates.Keys) {
    $vals = $forecastUpdates[$row]
    $wsForecast.Range("B$row").Value = "'" + $vals.B
    $wsForecast.Range("D$row").Value = $vals.D
}

# --- Sheet "Summary": refresh the computed statistics. All values in
# column B on this sheet are stored as text, including the numeric-looking
# ones, so every write uses a leading apostrophe to force text.

$summaryUpdates = @{
    2  = "2023-01-08 to 2025-01-05"
    4  = "102"
    6  = "42"
    8  = "4554 units"
    9  = "630"
    10 = "341"
    11 = "179"
    12 = "55"
    14 = "31"
    15 = "2025-04-27"
}

foreach ($row in $summaryUpdates.Keys) {
    $wsSummary.Range("B$row").Value = "'" + $summaryUpdates[$row]
}

Write-Output "Forecast rollover + summary refresh applied"
